$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying dataset (missing-data imputation benchmark) was regenerated with a
# different random "missingness" pattern:
#   - the "RM 232" record (old row 26) is no longer present at all
#   - the "SC 92" record (old row 28) is no longer present at all
#   - a handful of individual cells in column B ("A" measurement) flip which
#     record is missing its value
#
# First remove the two rows that disappeared entirely, which shifts every
# row below them up (mirrors what a human editor deleting worksheet rows
# would do).

# Delete old row 26 ("RM 232"); everything below shifts up by one.
$ws.Rows.Item(26).Delete()

# After that delete, the record that used to be "SC 92" (old row 28) is now
# at row 27; remove it too, shifting the remaining rows up again.
$ws.Rows.Item(27).Delete()

# Now rows 26-33 hold: SC 5, SC 101, SC 105, SC 119, SC 120, SC 132, SC 193, SC 232
# Column B ("A") missingness changed for three of these records - fix them up.

# SC 5 (row 26): previously missing in B, now has a value.
$ws.Range("B26").Value = -20.2

# SC 101 (row 27): previously had a value in B, now missing.
$ws.Range("B27").ClearContents()

# SC 232 (row 33): previously missing in B, now has a value.
$ws.Range("B33").Value = -19.5

# Keep the sheet's declared dimension in sync with the now-smaller used range.
$ws.Range("A1:F33").Select()
